$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D"  = 0.08220000000000001
    "G"  = 0.02665610979150172
    "H"  = 0.02665610979150172
    "I"  = -0.004011612562681447
    "J"  = -0.004011612562681447
    "K"  = -5.62
    "L"  = -0.0148324096067564
    "M"  = 3.42
    "N"  = 0.1425
    "O"  = -0.6085409252669038
    "P"  = 3.42
    "Q"  = 0.1425
    "R"  = -0.6085409252669038
    "U"  = 31.5
    "V"  = 1.3125
    "W"  = -0.03802435723951286
    "X"  = 0.3278012896401485
    "Y"  = -0.3658256468796613
    "Z"  = 2.228823529411764
    "AA" = -0.008941176470588235
    "AB" = 0.1307281099723083
    "AC" = -0.1396692864428966
    "AD" = 75.2
    "AE" = 0
    "AF" = 75.2
    "AG" = 43.7
    "AH" = 0.7580645161290323
    "AI" = 0.3369175627240144
    "AJ" = 0.6454948301329394
    "AK" = 0.2279603547209181
    "AL" = 3.27
    "AM" = 2.949
    "AN" = 4.39766081871345
    "AO" = -0.4648318042813456
    "AP" = 2.555555555555555
    "AQ" = -0.5154289589691421
}

foreach ($row in 2,3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
